$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 'Dr. Servinaz Sayed Mohammad, Dr. Veronia Rafat, Dr. Heba Mahmoud Ali, Dr. Amira Sobhy, Dr. Shimaa Ahmad Mekki, Dr. Rana Abo-Zaid, Dr. Hend Mahmoud, Dr. Alshimaa Atef'
$ws.Range("G3").Value = 'Dr. Gehan Adel, Administrator, Dr. Manar Montaser, Dr. Alshimaa Atef'
$ws.Range("G4").Value = 'Dr. Heba Mahmoud Ali, Dr. Menna tuâ€™Allah Medhat, Dr. Majorelle Magdy, Dr. Shimaa Ahmad Mekki, Dr. Nourhan Mahmoud, Dr. Asmaa Reda, Dr. Hanan Ragab'
$ws.Range("G6").Value = 'Dr. Sara Nabil, Dr. Safa Hany'
$ws.Range("G9").Value = 'Dr. Eman M. Abo-Sakaya, Dr. Yasmeena Fattoh, Dr. Marina Youhanna, Dr. Madeha Saeed'
$ws.Range("G12").Value = 'Dr. Mona Ibrahim Hussein, Dr. Dalia Tarek Elwan, Dr. Heba Al-Sayed Mohammad'
$ws.Range("G14").Value = 'Dr. Sarah Mahdy, Dr. Alaa Ashraf'
$ws.Range("G15").Value = 'Dr. Sarah Mahdy, Dr. Alaa Ashraf'
$ws.Range("G18").Value = 'Dr. Abdullah El-Agrody, Dr. Eman Samir Gabry, Dr. Ola Abd Al-Fattah, Dr. Wafaa Ebida'
$ws.Range("G19").Value = 'Dr. Marina Sorial, Dr. Neveen Nashaat, Dr. Yasmin, Dr. Eman Samir Gabry, Dr. Wafaa Ebida'
$ws.Range("G20").Value = 'Dr. Marina Sorial, Dr. Neveen Nashaat, Dr. Yasmin, Dr. Marina Atef, Dr. Nardine, Dr. Monica, Dr. Remon'
$ws.Range("G21").Value = 'Dr. Servinaz Sayed Mohammad, Dr. Veronia Rafat, Dr. Heba Mahmoud Ali, Dr. Amira Sobhy, Dr. Shimaa Ahmad Mekki, Dr. Rana Abo-Zaid, Dr. Hend Mahmoud, Dr. Alshimaa Atef'
$ws.Range("G22").Value = 'Dr. Gehan Adel, Administrator, Dr. Manar Montaser, Dr. Alshimaa Atef'
$ws.Range("G23").Value = 'Dr. Heba Mahmoud Ali, Dr. Menna tuâ€™Allah Medhat, Dr. Majorelle Magdy, Dr. Shimaa Ahmad Mekki, Dr. Nourhan Mahmoud, Dr. Asmaa Reda, Dr. Hanan Ragab'
$ws.Range("G24").Value = 'Dr. Fatma Elhady, Dr. Nada Mohammad, Dr. Lamiaa Ossama, Dr. Abeer Ragab, Dr. Amera Ahmad Saad'
$ws.Range("G28").Value = 'Dr. Yasmeena Fattoh, Dr. Dina Adel, Dr. Esraa Mostafa, Dr. Nourhan Osama, Dr. Marwa Mustafa, Dr. Arwa Al-Sayed, Dr. Basma Hamed, Dr. Eman M. Abo-Sakaya, Dr. Madeha Saeed, Dr. Sarah Abdelmohsen'
$ws.Range("G29").Value = 'Dr. Amira Ibrahim, Dr. Yasmeena Fattoh, Dr. Esraa Mostafa'
$ws.Range("G31").Value = 'Dr. Mona Ibrahim Hussein, Dr. Dalia Tarek Elwan, Dr. Heba Al-Sayed Mohammad'
$ws.Range("G34").Value = 'Dr. Sarah Mahdy, Dr. Alaa Ashraf'
$ws.Range("G37").Value = 'Dr. Abdullah El-Agrody, Dr. Eman Samir Gabry, Dr. Ola Abd Al-Fattah, Dr. Wafaa Ebida'
$ws.Range("G38").Value = 'Dr. Marina Sorial, Dr. Neveen Nashaat, Dr. Yasmin, Dr. Marina Atef, Dr. Nardine, Dr. Monica, Dr. Remon'
$ws.Range("G39").Value = 'Dr. Marina Sorial, Dr. Neveen Nashaat, Dr. Yasmin, Dr. Marina Atef, Dr. Nardine, Dr. Monica, Dr. Remon'
$ws.Range("G40").Value = 'Dr. Servinaz Sayed Mohammad, Dr. Veronia Rafat, Dr. Heba Mahmoud Ali, Dr. Amira Sobhy, Dr. Shimaa Ahmad Mekki, Dr. Rana Abo-Zaid, Dr. Hend Mahmoud, Dr. Alshimaa Atef'
$ws.Range("G41").Value = 'Dr. Amira Sobhy, Dr. Mohammad El-Tanany, Dr. Shimaa Ahmad Mekki, Dr. Hanan Ragab, Dr. Hend Mahmoud, Dr. Alshimaa Atef'
$ws.Range("G42").Value = 'Dr. Servinaz Sayed Mohammad, Dr. Menna tuâ€™Allah Medhat, Dr. Shimaa Ahmad Mekki, Dr. Eman Tantawi, Dr. Alshimaa Atef'
$ws.Range("G43").Value = 'Dr. Menna tu''Alllah Mohammad, Dr. Fatma Elhady, Dr. Kerelos Zareef, Dr. Nada Mohammad, Dr. Abeer Ragab, Dr. Lamiaa Ossama, Dr. Amera Ahmad Saad'
$ws.Range("G44").Value = 'Dr. Sara Nabil, Dr. Safa Hany'
$ws.Range("G47").Value = 'Dr. Esraa Mostafa, Dr. Nourhan Osama, Dr. Amira Ibrahim, Dr. Maryam Ahmad, Dr. Arwa Al-Sayed, Dr. Merna Said'
$ws.Range("G48").Value = 'Dr. Amany Raafat, Dr. Maryam Ahmad, Dr. Eman M. Abo-Sakaya, Dr. Yasmeena Fattoh, Dr. Sarah Abdelmohsen, Dr. Merna Said'
$ws.Range("G50").Value = 'Dr. Mona Ibrahim Hussein, Dr. Dalia Tarek Elwan, Dr. Heba Al-Sayed Mohammad'
$ws.Range("G52").Value = 'Dr. Sarah Mahdy, Dr. Alaa Ashraf'
$ws.Range("G54").Value = 'Dr. Afaf Abdallah, Dr. Amr Saeed'
$ws.Range("G56").Value = 'Dr. Abdullah El-Agrody, Dr. Eman Samir Gabry, Dr. Ola Abd Al-Fattah, Dr. Wafaa Ebida'
$ws.Range("G57").Value = 'Dr. Marina Sorial, Dr. Neveen Nashaat, Dr. Yasmin, Dr. Marina Atef, Dr. Nardine, Dr. Monica, Dr. Remon'
$ws.Range("G58").Value = 'Dr. Marina Sorial, Dr. Neveen Nashaat, Dr. Yasmin, Dr. Marina Atef, Dr. Nardine, Dr. Monica, Dr. Remon'
$ws.Range("G59").Value = 'Dr. Servinaz Sayed Mohammad, Dr. Heba Mahmoud Ali, Dr. Amira Sobhy, Dr. Mohammad El-Tanany, Dr. Nesma, Dr. Nourhan Mahmoud, Dr. Asmaa Reda'
$ws.Range("G60").Value = 'Dr. Amira Sobhy, Dr. Mohammad El-Tanany, Dr. Shimaa Ahmad Mekki, Dr. Hanan Ragab, Dr. Hend Mahmoud, Dr. Alshimaa Atef'
$ws.Range("G61").Value = 'Dr. Majorelle Magdy, Dr. Amira Sobhy, Dr. Shimaa Ahmad Mekki, Dr. Nahla Nagiub, Dr. Asmaa Reda'
$ws.Range("G66").Value = 'Dr. Madeha Saeed, Dr. Dina Adel, Dr. Amira Ibrahim, Dr. Eman M. Abo-Sakaya, Dr. Yasmeena Fattoh, Dr. Marina Youhanna'
$ws.Range("G67").Value = 'Dr. Amira Ibrahim, Dr. Yasmeena Fattoh, Dr. Esraa Mostafa'
$ws.Range("G71").Value = 'Dr. Nouran Mahmoud, Dr. Sarah Mahdy'
$ws.Range("G72").Value = 'Dr. Nouran Mahmoud, Dr. Sarah Mahdy'
$ws.Range("G75").Value = 'Dr. Abdullah El-Agrody, Dr. Eman Samir Gabry, Dr. Ola Abd Al-Fattah, Dr. Wafaa Ebida'
$ws.Range("G76").Value = 'Dr. Marina Sorial, Dr. Neveen Nashaat, Dr. Yasmin, Dr. Eman Samir Gabry, Dr. Wafaa Ebida'
$ws.Range("G77").Value = 'Dr. Marina Sorial, Dr. Neveen Nashaat, Dr. Yasmin, Dr. Marina Atef, Dr. Nardine, Dr. Monica, Dr. Remon'
$ws.Range("G78").Value = 'Dr. Servinaz Sayed Mohammad, Dr. Heba Mahmoud Ali, Dr. Amira Sobhy, Dr. Mohammad El-Tanany, Dr. Nesma, Dr. Nourhan Mahmoud, Dr. Asmaa Reda'
$ws.Range("G79").Value = 'Dr. Amira Sobhy, Dr. Mohammad El-Tanany, Dr. Shimaa Ahmad Mekki, Dr. Hanan Ragab, Dr. Hend Mahmoud, Dr. Alshimaa Atef'
$ws.Range("G80").Value = 'Dr. Majorelle Magdy, Dr. Amira Sobhy, Dr. Shimaa Ahmad Mekki, Dr. Nahla Nagiub, Dr. Asmaa Reda'
$ws.Range("G81").Value = 'Dr. Fatma Elhady, Dr. Nada Mohammad, Dr. Lamiaa Ossama, Dr. Abeer Ragab, Dr. Amera Ahmad Saad'
$ws.Range("G85").Value = 'Dr. Madeha Saeed, Dr. Dina Adel, Dr. Amira Ibrahim, Dr. Eman M. Abo-Sakaya, Dr. Yasmeena Fattoh, Dr. Marina Youhanna'
$ws.Range("G86").Value = 'Dr. Amany Raafat, Dr. Maryam Ahmad, Dr. Eman M. Abo-Sakaya, Dr. Yasmeena Fattoh, Dr. Sarah Abdelmohsen, Dr. Merna Said'
$ws.Range("G88").Value = 'Dr. Mona Ibrahim Hussein, Dr. Dalia Tarek Elwan, Dr. Heba Al-Sayed Mohammad'
$ws.Range("G90").Value = 'Dr. Nouran Mahmoud, Dr. Sarah Mahdy'
$ws.Range("G91").Value = 'Dr. Nouran Mahmoud, Dr. Sarah Mahdy'
$ws.Range("G94").Value = 'Dr. Abdullah El-Agrody, Dr. Eman Samir Gabry, Dr. Ola Abd Al-Fattah, Dr. Wafaa Ebida'
$ws.Range("G95").Value = 'Dr. Marina Sorial, Dr. Neveen Nashaat, Dr. Yasmin, Dr. Eman Samir Gabry, Dr. Wafaa Ebida'
$ws.Range("G96").Value = 'Dr. Marina Sorial, Dr. Neveen Nashaat, Dr. Yasmin, Dr. Marina Atef, Dr. Nardine, Dr. Monica, Dr. Remon'
$ws.Range("G97").Value = 'Dr. Servinaz Sayed Mohammad, Dr. Heba Mahmoud Ali, Dr. Amira Sobhy, Dr. Mohammad El-Tanany, Dr. Nesma, Dr. Nourhan Mahmoud, Dr. Asmaa Reda'
$ws.Range("G98").Value = 'Dr. Amira Sobhy, Dr. Mohammad El-Tanany, Dr. Shimaa Ahmad Mekki, Dr. Hanan Ragab, Dr. Hend Mahmoud, Dr. Alshimaa Atef'
$ws.Range("G99").Value = 'Dr. Servinaz Sayed Mohammad, Dr. Menna tuâ€™Allah Medhat, Dr. Shimaa Ahmad Mekki, Dr. Eman Tantawi, Dr. Alshimaa Atef'
$ws.Range("G100").Value = 'Dr. Menna tu''Alllah Mohammad, Dr. Fatma Elhady, Dr. Kerelos Zareef, Dr. Nada Mohammad, Dr. Abeer Ragab, Dr. Lamiaa Ossama, Dr. Amera Ahmad Saad'
$ws.Range("G104").Value = 'Dr. Esraa Mostafa, Dr. Nourhan Osama, Dr. Amira Ibrahim, Dr. Maryam Ahmad, Dr. Arwa Al-Sayed, Dr. Merna Said'
$ws.Range("G113").Value = 'Dr. Abdullah El-Agrody, Dr. Eman Samir Gabry, Dr. Ola Abd Al-Fattah, Dr. Wafaa Ebida'
